# Update on Base page
# - Application_Details!F2 (Modified Application Name) -> Atanu_Test_940
# - Step_Details!I2 (Modified Step Name) -> Atanu_Step_245

$wb = $excel.ActiveWorkbook

$wsApplicationDetails = $wb.Worksheets.Item("Application_Details")
$wsStepDetails = $wb.Worksheets.Item("Step_Details")

$wsApplicationDetails.Range("F2").Value = "Atanu_Test_940"
$wsStepDetails.Range("I2").Value = "Atanu_Step_245"
